$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "UsernameCredentials"

# --- Add the new sheet right after the first one ---
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "UserEmailCredentials"

# --- Update sheet1 (UsernameCredentials) header text ---
$ws1.Range("A1").Value2 = "User name"
$ws1.Range("B1").Value2 = "Password"

# --- Data for both sheets ---
$names = @("Vendetta", "Ashen", "Porto", "beautiful_Antonia", "Armstrong", "spaceX", "agulek", "wild_baby", "greasy-muffin", "Coookie")
$emails = @("Vendetta@test1.com", "Ashen@test2.com", "Porto@test3.com", "beautiful_Antonia@test4.com", "Armstrong@test5.com", "spaceX@test6.com", "agulek@test7.com", "wild_baby@test8.com", "greasy-muffin@test9.com", "Coookie@test10.com")
$passwords = @("pass1", "pass2", "pass3", "pass4", "pass5", "pass6", "pass7", "pass8", "pass9", "pass10")

# --- Sheet2 header ---
$ws2.Range("A1").Value2 = "User email"
$ws2.Range("B1").Value2 = "Password"
$ws2.Range("A1:B1").Font.Bold = $true

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2

    # sheet1: username + password (values already present, just re-affirm)
    $ws1.Range("A$row").Value2 = $names[$i]
    $ws1.Range("B$row").Value2 = $passwords[$i]

    # sheet2: email as hyperlink + password
    $mail = "mailto:" + $emails[$i]
    $ws2.Hyperlinks.Add($ws2.Range("A$row"), $mail, [System.Type]::Missing, [System.Type]::Missing, $emails[$i]) | Out-Null
    $ws2.Range("B$row").Value2 = $passwords[$i]
}

# --- Column widths ---
# Column A on both sheets ends up "best fit" to the same width (the longer
# e-mail values drive the sizing on both sheets), so apply the same width to
# both.
$ws1.Columns.Item(1).ColumnWidth = 27.43
$ws2.Columns.Item(1).ColumnWidth = 27.43
$ws2.Columns.Item(2).ColumnWidth = 12.15

# --- Selections ---
$ws1.Range("B8").Select() | Out-Null
$ws2.Range("D7").Select() | Out-Null

# --- Make sheet2 the active / displayed sheet ---
$ws2.Activate() | Out-Null
